$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cross review points ")

# Update the Decision column (G) status to "Accepted" for all review rows
$ws.Range("G2").Value = "Accepted"
$ws.Range("G3").Value = "Accepted"
$ws.Range("G4").Value = "Accepted"
$ws.Range("G5").Value = "Accepted"
$ws.Range("G6").Value = "Accepted"
$ws.Range("G7").Value = "Accepted"
$ws.Range("G8").Value = "Accepted"
$ws.Range("G9").Value = "Accepted"
$ws.Range("G10").Value = "Accepted"

$ws.Activate()
$ws.Range("G10").Select()
